$wb = $excel.ActiveWorkbook

# Append two new daily records (rows 33 and 34) to each of the four
# worksheets (FE_LFT_#1, FE_LFT_#2, FE_PLT_#1, FE_PLT_#2), continuing
# the existing time-series database with entries for 2025-06-11 and
# 2025-06-12. Column layout: A=time, B..E=raw hex fields (text),
# F..I=decoded decimal fields (numeric).

$dateFormat = "YYYY-MM-DD HH:MM:SS"

# --- Sheet 1: FE_LFT_#1 ---
$ws = $wb.Worksheets.Item(1)
$newRows = @(
  @{ r = 33; A = "45819.49663194444"; B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x01,0x70"; E = "0xf"; F = "380"; G = "7.598631275147109e+23"; H = "368"; I = "15" },
  @{ r = 34; A = "45820.49606481481"; B = "0x01,0x7c"; C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"; D = "0x01,0x70"; E = "0xf"; F = "380"; G = "7.598631275147109e+23"; H = "368"; I = "15" }
)
foreach ($row in $newRows) {
  $ws.Cells.Item($row.r, 1).Value = [double]$row.A
  $ws.Cells.Item($row.r, 1).NumberFormat = $dateFormat
  $ws.Cells.Item($row.r, 2).Value = $row.B
  $ws.Cells.Item($row.r, 3).Value = $row.C
  $ws.Cells.Item($row.r, 4).Value = $row.D
  $ws.Cells.Item($row.r, 5).Value = $row.E
  $ws.Cells.Item($row.r, 6).Value = [double]$row.F
  $ws.Cells.Item($row.r, 7).Value = [double]$row.G
  $ws.Cells.Item($row.r, 8).Value = [double]$row.H
  $ws.Cells.Item($row.r, 9).Value = [double]$row.I
}

# --- Sheet 2: FE_LFT_#2 ---
$ws = $wb.Worksheets.Item(2)
$newRows = @(
  @{ r = 33; A = "45819.49663194444"; B = "0x01,0x90"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x01,0x84"; E = "0xe"; F = "400"; G = "5.68432987514711e+23"; H = "388"; I = "14" },
  @{ r = 34; A = "45820.49606481481"; B = "0x01,0x90"; C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"; D = "0x01,0x80"; E = "0xe"; F = "400"; G = "5.68432987514711e+23"; H = "384"; I = "14" }
)
foreach ($row in $newRows) {
  $ws.Cells.Item($row.r, 1).Value = [double]$row.A
  $ws.Cells.Item($row.r, 1).NumberFormat = $dateFormat
  $ws.Cells.Item($row.r, 2).Value = $row.B
  $ws.Cells.Item($row.r, 3).Value = $row.C
  $ws.Cells.Item($row.r, 4).Value = $row.D
  $ws.Cells.Item($row.r, 5).Value = $row.E
  $ws.Cells.Item($row.r, 6).Value = [double]$row.F
  $ws.Cells.Item($row.r, 7).Value = [double]$row.G
  $ws.Cells.Item($row.r, 8).Value = [double]$row.H
  $ws.Cells.Item($row.r, 9).Value = [double]$row.I
}

# --- Sheet 3: FE_PLT_#1 ---
$ws = $wb.Worksheets.Item(3)
$newRows = @(
  @{ r = 33; A = "45819.49663194444"; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x6C"; E = "0x3"; F = "110"; G = "5.68631262647114e+23"; H = "108"; I = "3" },
  @{ r = 34; A = "45820.49606481481"; B = "0x00,0x6e"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; D = "0x00,0x6C"; E = "0x3"; F = "110"; G = "5.68631262647114e+23"; H = "108"; I = "3" }
)
foreach ($row in $newRows) {
  $ws.Cells.Item($row.r, 1).Value = [double]$row.A
  $ws.Cells.Item($row.r, 1).NumberFormat = $dateFormat
  $ws.Cells.Item($row.r, 2).Value = $row.B
  $ws.Cells.Item($row.r, 3).Value = $row.C
  $ws.Cells.Item($row.r, 4).Value = $row.D
  $ws.Cells.Item($row.r, 5).Value = $row.E
  $ws.Cells.Item($row.r, 6).Value = [double]$row.F
  $ws.Cells.Item($row.r, 7).Value = [double]$row.G
  $ws.Cells.Item($row.r, 8).Value = [double]$row.H
  $ws.Cells.Item($row.r, 9).Value = [double]$row.I
}

# --- Sheet 4: FE_PLT_#2 ---
$ws = $wb.Worksheets.Item(4)
$newRows = @(
  @{ r = 33; A = "45819.49663194444"; B = "0x00,0x6e"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x6C"; E = "0x3"; F = "110"; G = "9.85046333984776e+23"; H = "108"; I = "3" },
  @{ r = 34; A = "45820.49606481481"; B = "0x00,0x6e"; C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"; D = "0x00,0x6C"; E = "0x3"; F = "110"; G = "9.85046333984776e+23"; H = "108"; I = "3" }
)
foreach ($row in $newRows) {
  $ws.Cells.Item($row.r, 1).Value = [double]$row.A
  $ws.Cells.Item($row.r, 1).NumberFormat = $dateFormat
  $ws.Cells.Item($row.r, 2).Value = $row.B
  $ws.Cells.Item($row.r, 3).Value = $row.C
  $ws.Cells.Item($row.r, 4).Value = $row.D
  $ws.Cells.Item($row.r, 5).Value = $row.E
  $ws.Cells.Item($row.r, 6).Value = [double]$row.F
  $ws.Cells.Item($row.r, 7).Value = [double]$row.G
  $ws.Cells.Item($row.r, 8).Value = [double]$row.H
  $ws.Cells.Item($row.r, 9).Value = [double]$row.I
}
